$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.427.04"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.613.72"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.07"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.48"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.52"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("E10").Value = "  -1.61%  "

$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.373"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.077.19"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.68"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +4.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.424.89"
$ws.Range("D15").NumberFormat = "General"

$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.618.02"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.36"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("E19").Value = "  -1.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.79"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.92"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -3.27%  "

$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.532"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +1.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.48"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -1.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.161"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -1.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.00"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +3.47%  "

$ws.Range("E28").Value = "  +5.17%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.45"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +1.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.75"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +4.86%  "

$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.51"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("E34").Value = "  +8.57%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.29"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +4.36%  "

$ws.Range("E37").Value = "  +3.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "319.41"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +6.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.44"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +1.72%  "

$ws.Range("E40").Value = "  +2.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.851"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "136.34"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -2.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0994"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.52%  "

$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.95"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +1.32%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.610"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.68%  "

$ws.Range("E47").Value = "  +3.91%  "

$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("E49").Value = "  +1.32%  "

$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("E51").Value = "  +0.48%  "
